{"js": "// Apply the \"Juno: check in to OLPRODLOC\" localization/copy edits to the\n// Mystic Spice Premium Chai Tea product description.\n\nconst body = context.document.body;\n\nasync function replaceText(oldText, newText, opts) {\n  const options = Object.assign({ matchCase: true }, opts || {});\n  const results = body.search(oldText, options);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Title: brand name localized back to the English brand name.\nawait replaceText(\n  \": Mystisch Gew\u00fcrz Premium Chai Tee\",\n  \": Mystic Spice Premium Chai Tee\"\n);\n\n// \"Wichtige Merkmale:\" -> \"Wichtige Features:\" and make it bold (it was\n// missing the bold formatting the other section headings have).\n{\n  const results = body.search(\"Wichtige Merkmale:\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: Wichtige Merkmale:\");\n  }\n  const r = results.items[0];\n  r.font.bold = true;\n  r.insertText(\"Wichtige Features:\", \"Replace\");\n  await context.sync();\n}\n\nawait replaceText(\"Authentic Blend\", \"Authentische Mischung\");\n\nawait replaceText(\n  \": Unsere Chai ist eine harmonische Mischung aus Premium-Schwarzen Teebl\u00e4ttern und einer charakteristischen Auswahl an gemahlenen Gew\u00fcrzen, darunter Knoblauch, Karamom, Gerinnsel, Ingwer und schwarzer Pfeffer.\",\n  \": Unser Chai ist eine harmonische Mischung aus hochwertigen schwarzen Teebl\u00e4ttern und einer charakteristischen Auswahl an gemahlenen Gew\u00fcrzen wie Zimt, Kardamom, Nelken, Ingwer und schwarzem Pfeffer.\"\n);\n\nawait replaceText(\n  \"Gesundheitsf\u00f6rdernde Inhaltsstoffe: Jeder Bestandteil\",\n  \"Gesundheitsf\u00f6rdernde Inhaltsstoffe\"\n);\n\nawait replaceText(\n  \" von Mystisch Gew\u00fcrz-Chai-Tee wird f\u00fcr seine nat\u00fcrlichen Gesundheitlichen Vorteile ausgew\u00e4hlt.\",\n  \": Alle Inhaltsstoffe des Mystic Spice Chai Tea werden aufgrund ihrer nat\u00fcrlichen gesundheitsf\u00f6rdernden Eigenschaften ausgew\u00e4hlt.\"\n);\n\nawait replaceText(\n  \": Das warme, w\u00fcrzige Aroma und tiefe, belebende Geschmack unserer Chai machen es zum perfekten Getr\u00e4nk, um Ihren Tag zu beginnen oder sich am Abend zu entspannen.\",\n  \": Das warme, w\u00fcrzige Aroma und tiefe, belebende Geschmack unseres Chai machen ihn zum perfekten Getr\u00e4nk, um in den Tag zu starten oder am Abend zu entspannen.\"\n);\n\nawait replaceText(\"Vielseitige Brauoptionen\", \"Vielf\u00e4ltige Zubereitungsm\u00f6glichkeiten\");\n\nawait replaceText(\n  \": Ob Sie Ihre Chai hei\u00df dampfen, als erfrischender Eistee oder als cremefarbene Latte lieben, ist unsere Mischung vielseitig genug f\u00fcr jede Vorliebe.\",\n  \": Ob Sie Ihren Chai dampfend hei\u00df, als erfrischenden Eistee oder als cremigen Latte m\u00f6gen \u2013 unsere Mischung ist vielseitig genug, um allen Vorlieben gerecht zu werden.\"\n);\n\nawait replaceText(\n  \": Wir engagieren uns f\u00fcr Nachhaltigkeit, wir beziehen unsere Zutaten aus kleinfl\u00e4chigen Farmen, die \u00f6kologische Landwirtschaft betreiben, und sorgen nicht nur f\u00fcr die feinste Qualit\u00e4t, sondern auch f\u00fcr das Wohlergehen unseres Planeten.\",\n  \": Da wir uns der Nachhaltigkeit verpflichtet haben, beziehen wir unsere Zutaten von kleinen Bauernh\u00f6fen, die \u00f6kologische Landwirtschaft betreiben. So garantieren wir nicht nur beste Qualit\u00e4t, sondern tragen auch zum Wohlergehen unseres Planeten bei.\"\n);\n\nawait replaceText(\n  \": Mystisch Gew\u00fcrz-Chai-Tee kommt in wundersch\u00f6n gestalteten, umweltfreundlichen Verpackungen, sodass es ein ideales Geschenk f\u00fcr Teeliebhaber oder ein luxuri\u00f6ser Genuss f\u00fcr sich selbst ist.\",\n  \": Mystic Spice Chai Tea wird in einer wundersch\u00f6nen, umweltfreundlichen Verpackung geliefert, die ihn zu einem idealen Geschenk f\u00fcr Teeliebhaber oder zu einem luxuri\u00f6sen Genuss f\u00fcr Sie selbst macht.\"\n);\n\nawait replaceText(\n  \": Tee-Enthusiasten, gesundheitsbewusste Einzelpersonen, Liebhaber warmer, w\u00fcrziger Getr\u00e4nke und jeder, der die reichen Aromen der traditionellen indischen Chai erkunden m\u00f6chte.\",\n  \": Teeliebhaber, gesundheitsbewusste Menschen, Liebhaber von warmen, w\u00fcrzigen Getr\u00e4nken und alle, die den reichen Geschmack des traditionellen indischen Chai entdecken m\u00f6chten.\"\n);\n", "ps1": "# Apply the \"Juno: check in to OLPRODLOC\" localization/copy edits to the\n# Mystic Spice Premium Chai Tea product description.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $ok) {\n        Write-Output \"NOT FOUND: $old\"\n    }\n}\n\n# Title: brand name localized back to the English brand name.\nReplace-Text \": Mystisch Gew\u00fcrz Premium Chai Tee\" \": Mystic Spice Premium Chai Tee\"\n\n# \"Wichtige Merkmale:\" -> \"Wichtige Features:\" and make it bold (it was\n# missing the bold formatting the other section headings have).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Replacement.Font.Bold = 1\n$find.Execute(\"Wichtige Merkmale:\", $false, $false, $false, $false, $false, $true, 1, $true, \"Wichtige Features:\", 2) | Out-Null\n\nReplace-Text \"Authentic Blend\" \"Authentische Mischung\"\n\nReplace-Text \": Unsere Chai ist eine harmonische Mischung aus Premium-Schwarzen Teebl\u00e4ttern und einer charakteristischen Auswahl an gemahlenen Gew\u00fcrzen, darunter Knoblauch, Karamom, Gerinnsel, Ingwer und schwarzer Pfeffer.\" \": Unser Chai ist eine harmonische Mischung aus hochwertigen schwarzen Teebl\u00e4ttern und einer charakteristischen Auswahl an gemahlenen Gew\u00fcrzen wie Zimt, Kardamom, Nelken, Ingwer und schwarzem Pfeffer.\"\n\nReplace-Text \"Gesundheitsf\u00f6rdernde Inhaltsstoffe: Jeder Bestandteil\" \"Gesundheitsf\u00f6rdernde Inhaltsstoffe\"\n\nReplace-Text \" von Mystisch Gew\u00fcrz-Chai-Tee wird f\u00fcr seine nat\u00fcrlichen Gesundheitlichen Vorteile ausgew\u00e4hlt.\" \": Alle Inhaltsstoffe des Mystic Spice Chai Tea werden aufgrund ihrer nat\u00fcrlichen gesundheitsf\u00f6rdernden Eigenschaften ausgew\u00e4hlt.\"\n\nReplace-Text \": Das warme, w\u00fcrzige Aroma und tiefe, belebende Geschmack unserer Chai machen es zum perfekten Getr\u00e4nk, um Ihren Tag zu beginnen oder sich am Abend zu entspannen.\" \": Das warme, w\u00fcrzige Aroma und tiefe, belebende Geschmack unseres Chai machen ihn zum perfekten Getr\u00e4nk, um in den Tag zu starten oder am Abend zu entspannen.\"\n\nReplace-Text \"Vielseitige Brauoptionen\" \"Vielf\u00e4ltige Zubereitungsm\u00f6glichkeiten\"\n\nReplace-Text \": Ob Sie Ihre Chai hei\u00df dampfen, als erfrischender Eistee oder als cremefarbene Latte lieben, ist unsere Mischung vielseitig genug f\u00fcr jede Vorliebe.\" \": Ob Sie Ihren Chai dampfend hei\u00df, als erfrischenden Eistee oder als cremigen Latte m\u00f6gen \u2013 unsere Mischung ist vielseitig genug, um allen Vorlieben gerecht zu werden.\"\n\nReplace-Text \": Wir engagieren uns f\u00fcr Nachhaltigkeit, wir beziehen unsere Zutaten aus kleinfl\u00e4chigen Farmen, die \u00f6kologische Landwirtschaft betreiben, und sorgen nicht nur f\u00fcr die feinste Qualit\u00e4t, sondern auch f\u00fcr das Wohlergehen unseres Planeten.\" \": Da wir uns der Nachhaltigkeit verpflichtet haben, beziehen wir unsere Zutaten von kleinen Bauernh\u00f6fen, die \u00f6kologische Landwirtschaft betreiben. So garantieren wir nicht nur beste Qualit\u00e4t, sondern tragen auch zum Wohlergehen unseres Planeten bei.\"\n\nReplace-Text \": Mystisch Gew\u00fcrz-Chai-Tee kommt in wundersch\u00f6n gestalteten, umweltfreundlichen Verpackungen, sodass es ein ideales Geschenk f\u00fcr Teeliebhaber oder ein luxuri\u00f6ser Genuss f\u00fcr sich selbst ist.\" \": Mystic Spice Chai Tea wird in einer wundersch\u00f6nen, umweltfreundlichen Verpackung geliefert, die ihn zu einem idealen Geschenk f\u00fcr Teeliebhaber oder zu einem luxuri\u00f6sen Genuss f\u00fcr Sie selbst macht.\"\n\nReplace-Text \": Tee-Enthusiasten, gesundheitsbewusste Einzelpersonen, Liebhaber warmer, w\u00fcrziger Getr\u00e4nke und jeder, der die reichen Aromen der traditionellen indischen Chai erkunden m\u00f6chte.\" \": Teeliebhaber, gesundheitsbewusste Menschen, Liebhaber von warmen, w\u00fcrzigen Getr\u00e4nken und alle, die den reichen Geschmack des traditionellen indischen Chai entdecken m\u00f6chten.\"\n"}
